# Apply cryptocurrency price/volume updates (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) store plain text that looks numeric
# (e.g. "494.24", "0.0000137", "  +1.32%  "). Force Text format first so
# Excel keeps the exact literal string instead of coercing it to a Double
# (which would also silently drop meaningful trailing zeros).
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "56.592.88"
$ws.Range("E2").Value = "  +1.32%  "

$ws.Range("D3").Value = "2.494.11"
$ws.Range("E3").Value = "  -0.99%  "

$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").Value = "494.24"
$ws.Range("E5").Value = "  +1.71%  "

$ws.Range("D6").Value = "153.04"
$ws.Range("E6").Value = "  +7.27%  "

$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("D8").Value = "0.512"
$ws.Range("E8").Value = "  -0.29%  "

$ws.Range("D9").Value = "2.507.10"
$ws.Range("E9").Value = "  -0.26%  "

$ws.Range("D10").Value = "5.78"
$ws.Range("E10").Value = "  +5.00%  "

$ws.Range("D11").Value = "0.0989"
$ws.Range("E11").Value = "  -0.55%  "

$ws.Range("D12").Value = "0.335"
$ws.Range("E12").Value = "  +1.34%  "

$ws.Range("E13").Value = "  +1.22%  "

$ws.Range("D14").Value = "2.933.20"
$ws.Range("E14").Value = "  -0.49%  "

$ws.Range("D15").Value = "56.804.31"
$ws.Range("E15").Value = "  +1.72%  "

$ws.Range("D16").Value = "21.43"
$ws.Range("E16").Value = "  +3.43%  "

$ws.Range("D17").Value = "0.0000137"
$ws.Range("E17").Value = "  -1.24%  "

$ws.Range("D18").Value = "2.511.52"
$ws.Range("E18").Value = "  -0.12%  "

$ws.Range("D19").Value = "4.55"
$ws.Range("E19").Value = "  +3.34%  "

$ws.Range("D20").Value = "10.34"
$ws.Range("E20").Value = "  +2.52%  "

$ws.Range("D21").Value = "323.01"
$ws.Range("E21").Value = "  +0.37%  "

$ws.Range("D22").Value = "0.996"
$ws.Range("E22").Value = "  -0.26%  "

$ws.Range("D23").Value = "5.90"
$ws.Range("E23").Value = "  +2.92%  "

$ws.Range("D24").Value = "58.99"
$ws.Range("E24").Value = "  +1.36%  "

$ws.Range("D25").Value = "0.411"
$ws.Range("E25").Value = "  +0.45%  "

$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.09%  "

$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").Value = "0.163"
$ws.Range("E27").Value = "  -2.78%  "

$ws.Range("D28").Value = "2.616.08"
$ws.Range("E28").Value = "  +0.15%  "

$ws.Range("D29").Value = "7.61"
$ws.Range("E29").Value = "  +1.84%  "

$ws.Range("D30").Value = "0.0₃0815"
$ws.Range("E30").Value = "  +0.64%  "

$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  -0.11%  "

$ws.Range("D32").Value = "152.31"
$ws.Range("E32").Value = "  +1.75%  "

$ws.Range("D33").Value = "18.38"
$ws.Range("E33").Value = "  +0.74%  "

$ws.Range("D34").Value = "1.53"
$ws.Range("E34").Value = "  +2.25%  "

$ws.Range("D35").Value = "5.26"
$ws.Range("E35").Value = "  +0.83%  "

$ws.Range("E36").Value = "  +3.16%  "

$ws.Range("D37").Value = "3.80"
$ws.Range("E37").Value = "  +1.59%  "

$ws.Range("D38").Value = "0.868"
$ws.Range("E38").Value = "  -1.53%  "

$ws.Range("E39").Value = "  +4.30%  "

$ws.Range("D40").Value = "34.11"
$ws.Range("E40").Value = "  -0.53%  "

$ws.Range("D41").Value = "3.52"
$ws.Range("E41").Value = "  +2.50%  "

$ws.Range("D42").Value = "0.0565"
$ws.Range("E42").Value = "  +1.73%  "

$ws.Range("D43").Value = "0.618"
$ws.Range("E43").Value = "  +0.26%  "

$ws.Range("D44").Value = "0.995"
$ws.Range("E44").Value = "  -0.24%  "

$ws.Range("D45").Value = "4.94"
$ws.Range("E45").Value = "  +5.18%  "

$ws.Range("D46").Value = "268.81"
$ws.Range("E46").Value = "  +5.37%  "

$ws.Range("D47").Value = "0.0930"
$ws.Range("E47").Value = "  +2.41%  "

$ws.Range("D48").Value = "0.0231"
$ws.Range("E48").Value = "  +2.41%  "

$ws.Range("D49").Value = "10.21"
$ws.Range("E49").Value = "  +0.40%  "

$ws.Range("D50").Value = "17.89"
$ws.Range("E50").Value = "  +1.08%  "

$ws.Range("D51").Value = "1.905.84"
$ws.Range("E51").Value = "  -4.73%  "

# Restore the default cell style (the sheet never used an explicit Text
# style on these cells) now that the literal values are locked in.
$dataRange.Style = "Normal"
